# Listas sem duplicação de professores
# Cleans up cells that showed the teacher's own class duplicated across
# all weekday slots; keeps only the correct single slot (or a plain "-"
# when no slot remains), matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Wednesday (D) loses the duplicated list, becomes a plain dash
$ws.Range("D3").Value = "-"

# Row 4 - Wednesday (D) becomes a plain dash, Thursday (E) gains the
# reshuffled list (teacher's slot now in 3rd position)
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "[-, -, 'MCT-3A-Eletropneumática', -]"

# Row 6 - same pattern as row 4
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "[-, -, 'MCT-3A-Eletropneumática', -]"

# Row 7 - same pattern as row 4
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "[-, -, 'MCT-3A-Eletropneumática', -]"

# Row 8 - Thursday (E) gains the reshuffled list
$ws.Range("E8").Value = "[-, -, 'MCT-3A-Eletropneumática', -]"

# Rows 11, 12, 14, 15 - Friday (F) loses the duplicated list
$ws.Range("F11").Value = "-"
$ws.Range("F12").Value = "-"
$ws.Range("F14").Value = "-"
$ws.Range("F15").Value = "-"

# Rows 18-21 - Monday (B) and Tuesday (C) lose duplicated lists; Friday
# (F) loses the redundant single entry in rows 18 and 19
$ws.Range("B18").Value = "-"
$ws.Range("C18").Value = "-"
$ws.Range("F18").Value = "-"

$ws.Range("B19").Value = "-"
$ws.Range("C19").Value = "-"
$ws.Range("F19").Value = "-"

$ws.Range("B20").Value = "-"
$ws.Range("C20").Value = "-"

$ws.Range("B21").Value = "-"
